$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

# Update the scraped_at timestamps (column K) for rows 2-48 to reflect
# the latest scrape run.
$ws.Range("K2").Value = "2025-11-03T11:34:29.777784+00:00"
$ws.Range("K3").Value = "2025-11-03T11:34:32.305167+00:00"
$ws.Range("K4").Value = "2025-11-03T11:34:32.305220+00:00"
$ws.Range("K5").Value = "2025-11-03T11:34:32.305229+00:00"
$ws.Range("K6").Value = "2025-11-03T11:34:32.305236+00:00"
$ws.Range("K7").Value = "2025-11-03T11:34:34.306472+00:00"
$ws.Range("K8").Value = "2025-11-03T11:34:34.306489+00:00"
$ws.Range("K9").Value = "2025-11-03T11:34:34.306498+00:00"
$ws.Range("K10").Value = "2025-11-03T11:34:36.390291+00:00"
$ws.Range("K11").Value = "2025-11-03T11:34:38.883808+00:00"
$ws.Range("K12").Value = "2025-11-03T11:34:38.883826+00:00"
$ws.Range("K13").Value = "2025-11-03T11:34:38.883835+00:00"
$ws.Range("K14").Value = "2025-11-03T11:34:40.942724+00:00"
$ws.Range("K15").Value = "2025-11-03T11:34:40.942739+00:00"
$ws.Range("K16").Value = "2025-11-03T11:34:40.942747+00:00"
$ws.Range("K17").Value = "2025-11-03T11:34:47.417874+00:00"
$ws.Range("K18").Value = "2025-11-03T11:34:49.392150+00:00"
$ws.Range("K19").Value = "2025-11-03T11:34:51.403987+00:00"
$ws.Range("K20").Value = "2025-11-03T11:34:53.417579+00:00"
$ws.Range("K21").Value = "2025-11-03T11:34:53.417626+00:00"
$ws.Range("K22").Value = "2025-11-03T11:34:53.417650+00:00"
$ws.Range("K23").Value = "2025-11-03T11:34:55.436531+00:00"
$ws.Range("K24").Value = "2025-11-03T11:34:55.436564+00:00"
$ws.Range("K25").Value = "2025-11-03T11:34:55.436584+00:00"
$ws.Range("K26").Value = "2025-11-03T11:34:55.436603+00:00"
$ws.Range("K27").Value = "2025-11-03T11:35:01.960182+00:00"
$ws.Range("K28").Value = "2025-11-03T11:35:01.960211+00:00"
$ws.Range("K29").Value = "2025-11-03T11:35:01.960230+00:00"
$ws.Range("K30").Value = "2025-11-03T11:35:01.960247+00:00"
$ws.Range("K31").Value = "2025-11-03T11:35:01.960263+00:00"
$ws.Range("K32").Value = "2025-11-03T11:35:04.845955+00:00"
$ws.Range("K33").Value = "2025-11-03T11:35:04.845972+00:00"
$ws.Range("K34").Value = "2025-11-03T11:35:04.845980+00:00"
$ws.Range("K35").Value = "2025-11-03T11:35:06.879023+00:00"
$ws.Range("K36").Value = "2025-11-03T11:35:06.879039+00:00"
$ws.Range("K37").Value = "2025-11-03T11:35:06.879047+00:00"
$ws.Range("K38").Value = "2025-11-03T11:35:06.879057+00:00"
$ws.Range("K39").Value = "2025-11-03T11:35:06.879064+00:00"
$ws.Range("K40").Value = "2025-11-03T11:35:06.879071+00:00"
$ws.Range("K41").Value = "2025-11-03T11:35:06.879078+00:00"
$ws.Range("K42").Value = "2025-11-03T11:35:06.879084+00:00"
$ws.Range("K43").Value = "2025-11-03T11:35:09.438651+00:00"
$ws.Range("K44").Value = "2025-11-03T11:35:09.438671+00:00"
$ws.Range("K45").Value = "2025-11-03T11:35:14.477756+00:00"
$ws.Range("K46").Value = "2025-11-03T11:35:16.558785+00:00"
$ws.Range("K47").Value = "2025-11-03T11:35:16.558804+00:00"
$ws.Range("K48").Value = "2025-11-03T11:35:16.558812+00:00"
